$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (tab name)
$ws.Name = "Through 2022-10-14"

# Update header label in I1 (shared string "2022 (through 10-13)" -> "2022 (through 10-14)")
$ws.Range("I1").Value = "2022 (through 10-14)"

# Update data values
$ws.Range("I11").Value = 45
$ws.Range("I14").Value = 1323
